# Auto-generated edit script: updates cached market-price / profit
# columns (H:N) on specific rows across multiple crafting-job sheets,
# reflecting a refreshed pull from the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 5874.75
$ws.Range("I64").Value = 6233
$ws.Range("J64").Value = 4800
$ws.Range("K64").Value = 6233
$ws.Range("L64").Value = 4800
$ws.Range("M64").Value = -5985
$ws.Range("N64").Value = -5296

# Row 67
$ws.Range("H67").Value = 5874.75
$ws.Range("I67").Value = 6233
$ws.Range("J67").Value = 4800
$ws.Range("K67").Value = 6233
$ws.Range("L67").Value = 4800
$ws.Range("M67").Value = -5375
$ws.Range("N67").Value = -6516

# Row 98
$ws.Range("H98").Value = 2153.423
$ws.Range("I98").Value = 1787.45
$ws.Range("J98").Value = 3373.3333
$ws.Range("K98").Value = 1787.45
$ws.Range("L98").Value = 3373.3333
$ws.Range("M98").Value = -289.45
$ws.Range("N98").Value = -6369.3333

# Row 116
$ws.Range("H116").Value = 4171.357
$ws.Range("I116").Value = 5599.8
$ws.Range("J116").Value = 3377.7778
$ws.Range("K116").Value = 5599.8
$ws.Range("L116").Value = 3377.7778
$ws.Range("M116").Value = -2157.8
$ws.Range("N116").Value = -10261.7778

# Row 122
$ws.Range("H122").Value = 2153.423
$ws.Range("I122").Value = 1787.45
$ws.Range("J122").Value = 3373.3333
$ws.Range("K122").Value = 5362.35
$ws.Range("L122").Value = 10119.9999
$ws.Range("M122").Value = -2912.35
$ws.Range("N122").Value = -15019.9999

# Row 129
$ws.Range("H129").Value = 904361.5
$ws.Range("I129").Value = 264.30768
$ws.Range("J129").Value = 1324120.9
$ws.Range("K129").Value = 792.92304
$ws.Range("L129").Value = 3972362.7
$ws.Range("M129").Value = 4207.07696
$ws.Range("N129").Value = -3982362.7

# Row 138
$ws.Range("H138").Value = 1465.58
$ws.Range("I138").Value = 705.1356
$ws.Range("J138").Value = 2559.878
$ws.Range("K138").Value = 2115.4068
$ws.Range("L138").Value = 7679.634
$ws.Range("M138").Value = 3024.5932
$ws.Range("N138").Value = -17959.634


# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 17740.572
$ws.Range("I32").Value = 20327.059
$ws.Range("J32").Value = 9628.409
$ws.Range("K32").Value = 20327.059
$ws.Range("L32").Value = 9628.409
$ws.Range("M32").Value = -20040.059
$ws.Range("N32").Value = -10202.409


# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 150.6
$ws.Range("I5").Value = 150.6
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 150.6
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -37.59999999999999

# Row 134
$ws.Range("H134").Value = 15790.743
$ws.Range("I134").Value = 1171.3606
$ws.Range("J134").Value = 114877.664
$ws.Range("K134").Value = 3514.0818
$ws.Range("L134").Value = 344632.992
$ws.Range("M134").Value = -979.0817999999999
$ws.Range("N134").Value = -349702.992


# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()

# Row 31
$ws.Range("H31").Value = 1907.8108
$ws.Range("I31").Value = 1767.4822
$ws.Range("J31").Value = 2344.389
$ws.Range("K31").Value = 1767.4822
$ws.Range("L31").Value = 2344.389
$ws.Range("M31").Value = -1472.4822
$ws.Range("N31").Value = -2934.389

# Row 34
$ws.Range("H34").Value = 1907.8108
$ws.Range("I34").Value = 1767.4822
$ws.Range("J34").Value = 2344.389
$ws.Range("K34").Value = 1767.4822
$ws.Range("L34").Value = 2344.389
$ws.Range("M34").Value = -1565.4822
$ws.Range("N34").Value = -2748.389

# Row 62
$ws.Range("H62").Value = 8593.333000000001
$ws.Range("I62").Value = 9400
$ws.Range("J62").Value = 6980
$ws.Range("K62").Value = 9400
$ws.Range("L62").Value = 6980
$ws.Range("M62").Value = -8776
$ws.Range("N62").Value = -8228

# Row 65
$ws.Range("H65").Value = 8593.333000000001
$ws.Range("I65").Value = 9400
$ws.Range("J65").Value = 6980
$ws.Range("K65").Value = 47000
$ws.Range("L65").Value = 34900
$ws.Range("M65").Value = -43880
$ws.Range("N65").Value = -41140

# Row 94
$ws.Range("H94").Value = 2358.081
$ws.Range("I94").Value = 2754.889
$ws.Range("J94").Value = 2230.5356
$ws.Range("K94").Value = 2754.889
$ws.Range("L94").Value = 2230.5356
$ws.Range("M94").Value = -2303.889
$ws.Range("N94").Value = -3132.5356

# Row 132
$ws.Range("H132").Value = 1508.5193
$ws.Range("I132").Value = 1059.3158
$ws.Range("J132").Value = 2727.7856
$ws.Range("K132").Value = 3177.9474
$ws.Range("L132").Value = 8183.3568
$ws.Range("M132").Value = -647.9474
$ws.Range("N132").Value = -13243.3568

# Row 134
$ws.Range("H134").Value = 1269.8334
$ws.Range("I134").Value = 1141.42
$ws.Range("J134").Value = 2875
$ws.Range("K134").Value = 3424.26
$ws.Range("L134").Value = 8625
$ws.Range("M134").Value = -889.2600000000002
$ws.Range("N134").Value = -13695


# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 723461.0600000001
$ws.Range("I4").Value = 5050530
$ws.Range("J4").Value = 2282.9167
$ws.Range("K4").Value = 15151590
$ws.Range("L4").Value = 6848.750100000001
$ws.Range("M4").Value = -15151478
$ws.Range("N4").Value = -7072.750100000001

# Row 5
$ws.Range("H5").Value = 1468.2142
$ws.Range("I5").Value = 1828.8889
$ws.Range("J5").Value = 819
$ws.Range("K5").Value = 5486.6667
$ws.Range("L5").Value = 2457
$ws.Range("M5").Value = -5374.6667
$ws.Range("N5").Value = -2681

# Row 75
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()

# Row 78
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()

# Row 121
$ws.Range("H121").Value = 996.46155
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 996.46155
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 2989.38465
$ws.Range("N121").Value = -5609.38465

# Row 122
$ws.Range("H122").Value = 1251158.8
$ws.Range("I122").Value = 800
$ws.Range("J122").Value = 1429781.4
$ws.Range("K122").Value = 7200
$ws.Range("L122").Value = 12868032.6
$ws.Range("M122").Value = -4750
$ws.Range("N122").Value = -12872932.6

# Row 135
$ws.Range("H135").Value = 1468.2142
$ws.Range("I135").Value = 1828.8889
$ws.Range("J135").Value = 819
$ws.Range("K135").Value = 16460.0001
$ws.Range("L135").Value = 7371
$ws.Range("M135").Value = -13925.0001
$ws.Range("N135").Value = -12441


# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4294.355
$ws.Range("I70").Value = 4091.6365
$ws.Range("J70").Value = 4789.8887
$ws.Range("K70").Value = 4091.6365
$ws.Range("L70").Value = 4789.8887
$ws.Range("M70").Value = -3821.6365
$ws.Range("N70").Value = -5329.8887

# Row 73
$ws.Range("H73").Value = 4294.355
$ws.Range("I73").Value = 4091.6365
$ws.Range("J73").Value = 4789.8887
$ws.Range("K73").Value = 4091.6365
$ws.Range("L73").Value = 4789.8887
$ws.Range("M73").Value = -3155.6365
$ws.Range("N73").Value = -6661.8887

# Row 122
$ws.Range("H122").Value = 732366.5600000001
$ws.Range("I122").Value = 1317129
$ws.Range("J122").Value = 1413.5
$ws.Range("K122").Value = 3951387
$ws.Range("L122").Value = 4240.5
$ws.Range("M122").Value = -3948937
$ws.Range("N122").Value = -9140.5


# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 920265
$ws.Range("I40").Value = 1264576.2
$ws.Range("J40").Value = 2101.6667
$ws.Range("K40").Value = 1264576.2
$ws.Range("L40").Value = 2101.6667
$ws.Range("M40").Value = -1264440.2
$ws.Range("N40").Value = -2373.6667

# Row 132
$ws.Range("H132").Value = 1428.1642
$ws.Range("I132").Value = 1233.0167
$ws.Range("J132").Value = 3100.8572
$ws.Range("K132").Value = 3699.050099999999
$ws.Range("L132").Value = 9302.571599999999
$ws.Range("M132").Value = -1169.050099999999
$ws.Range("N132").Value = -14362.5716

# Row 136
$ws.Range("H136").Value = 3621.6667
$ws.Range("I136").Value = 1845
$ws.Range("J136").Value = 20500
$ws.Range("K136").Value = 5535
$ws.Range("L136").Value = 61500
$ws.Range("M136").Value = -2985
$ws.Range("N136").Value = -66600


# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 74
$ws.Range("H74").Value = 8716.5
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 8716.5
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 8716.5
$ws.Range("N74").Value = -10588.5
$ws.Range("M74").ClearContents()

# Row 77
$ws.Range("H77").Value = 8716.5
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 8716.5
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 26149.5
$ws.Range("N77").Value = -35509.5
$ws.Range("M77").ClearContents()

# Row 126
$ws.Range("H126").Value = 1041
$ws.Range("I126").Value = 1031.1666
$ws.Range("J126").Value = 1100
$ws.Range("K126").Value = 3093.4998
$ws.Range("L126").Value = 3300
$ws.Range("M126").Value = -623.4998000000001
$ws.Range("N126").Value = -8240

# Row 132
$ws.Range("H132").Value = 1230.7778
$ws.Range("I132").Value = 962.8077
$ws.Range("J132").Value = 1927.5
$ws.Range("K132").Value = 2888.4231
$ws.Range("L132").Value = 5782.5
$ws.Range("M132").Value = -358.4231
$ws.Range("N132").Value = -10842.5

# Row 136
$ws.Range("H136").Value = 553.1627999999999
$ws.Range("I136").Value = 247
$ws.Range("J136").Value = 2880
$ws.Range("K136").Value = 741
$ws.Range("L136").Value = 8640
$ws.Range("M136").Value = 1809
$ws.Range("N136").Value = -13740

